$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# G3 changes from 0 to 1
$ws.Range("G3").Value = 1

# H3:H18 change from 0 to 1
$ws.Range("H3:H18").Value = 1
